$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 25854
$ws.Range("J63").Value = 25854
$ws.Range("L63").Value = 25854
$ws.Range("N63").Value = -27102
$ws.Range("H66").Value = 25854
$ws.Range("J66").Value = 25854
$ws.Range("L66").Value = 77562
$ws.Range("N66").Value = -83802
$ws.Range("H100").Value = 4876.08
$ws.Range("I100").Value = 2357
$ws.Range("J100").Value = 6555.467
$ws.Range("K100").Value = 2357
$ws.Range("L100").Value = 6555.467
$ws.Range("M100").Value = -1816
$ws.Range("N100").Value = -7637.467
$ws.Range("H116").Value = 161606.53
$ws.Range("I116").Value = 208888.5
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 208888.5
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = -205446.5
$ws.Range("N116").Value = -10884

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1403.6522
$ws.Range("I45").Value = 1205.6
$ws.Range("J45").Value = 1775
$ws.Range("K45").Value = 1205.6
$ws.Range("L45").Value = 1775
$ws.Range("M45").Value = -828.5999999999999
$ws.Range("N45").Value = -2529
$ws.Range("H97").Value = 572.4186
$ws.Range("I97").Value = 548.7692
$ws.Range("J97").Value = 803
$ws.Range("K97").Value = 548.7692
$ws.Range("L97").Value = 803
$ws.Range("M97").Value = -52.76919999999996
$ws.Range("N97").Value = -1795
$ws.Range("H102").Value = 4581.143
$ws.Range("I102").Value = 1878.0714
$ws.Range("J102").Value = 9987.286
$ws.Range("K102").Value = 1878.0714
$ws.Range("L102").Value = 9987.286
$ws.Range("M102").Value = -256.0714
$ws.Range("N102").Value = -13231.286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""
$ws.Range("H80").Value = 391.35715
$ws.Range("I80").Value = 135
$ws.Range("J80").Value = 533.7778
$ws.Range("K80").Value = 135
$ws.Range("L80").Value = 533.7778
$ws.Range("M80").Value = 863
$ws.Range("N80").Value = -2529.7778
$ws.Range("H83").Value = 391.35715
$ws.Range("I83").Value = 135
$ws.Range("J83").Value = 533.7778
$ws.Range("K83").Value = 675
$ws.Range("L83").Value = 2668.889
$ws.Range("M83").Value = 4317
$ws.Range("N83").Value = -12652.889
$ws.Range("H86").Value = 7237.684
$ws.Range("I86").Value = 6236.909
$ws.Range("K86").Value = 6236.909
$ws.Range("M86").Value = -5113.909
$ws.Range("H89").Value = 7237.684
$ws.Range("I89").Value = 6236.909
$ws.Range("K89").Value = 31184.545
$ws.Range("M89").Value = -25568.545
$ws.Range("H94").Value = 819.48486
$ws.Range("I94").Value = 614.375
$ws.Range("J94").Value = 1366.4445
$ws.Range("K94").Value = 614.375
$ws.Range("L94").Value = 1366.4445
$ws.Range("M94").Value = -163.375
$ws.Range("N94").Value = -2268.4445
$ws.Range("H99").Value = 2973.348
$ws.Range("I99").Value = 3269.4
$ws.Range("J99").Value = 999.6667
$ws.Range("K99").Value = 3269.4
$ws.Range("L99").Value = 999.6667
$ws.Range("M99").Value = -1771.4
$ws.Range("N99").Value = -3995.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1321.6595
$ws.Range("I58").Value = 836.0303
$ws.Range("J58").Value = 2466.3572
$ws.Range("K58").Value = 836.0303
$ws.Range("L58").Value = 2466.3572
$ws.Range("M58").Value = -633.0303
$ws.Range("N58").Value = -2872.3572
$ws.Range("H86").Value = 5762
$ws.Range("I86").Value = 3014.4
$ws.Range("J86").Value = 19500
$ws.Range("K86").Value = 3014.4
$ws.Range("L86").Value = 19500
$ws.Range("M86").Value = -1891.4
$ws.Range("N86").Value = -21746
$ws.Range("H89").Value = 5762
$ws.Range("I89").Value = 3014.4
$ws.Range("J89").Value = 19500
$ws.Range("K89").Value = 15072
$ws.Range("L89").Value = 97500
$ws.Range("M89").Value = -9456
$ws.Range("N89").Value = -108732
$ws.Range("H107").Value = 301.46155
$ws.Range("I107").Value = 159.91667
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 159.91667
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 1760.08333
$ws.Range("N107").Value = -5840
$ws.Range("H136").Value = 1321.6595
$ws.Range("I136").Value = 836.0303
$ws.Range("J136").Value = 2466.3572
$ws.Range("K136").Value = 2508.0909
$ws.Range("L136").Value = 7399.071599999999
$ws.Range("M136").Value = 41.90909999999985
$ws.Range("N136").Value = -12499.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 41.761906
$ws.Range("J12").Value = 53
$ws.Range("L12").Value = 159
$ws.Range("N12").Value = -505
$ws.Range("H33").Value = 163.2
$ws.Range("I33").Value = 30
$ws.Range("J33").Value = 196.5
$ws.Range("K33").Value = 180
$ws.Range("L33").Value = 1179
$ws.Range("M33").Value = 103
$ws.Range("N33").Value = -1745
$ws.Range("H131").Value = 884.40845
$ws.Range("I131").Value = 477.5
$ws.Range("J131").Value = 908.7015
$ws.Range("K131").Value = 1432.5
$ws.Range("L131").Value = 2726.1045
$ws.Range("M131").Value = 3607.5
$ws.Range("N131").Value = -12806.1045
$ws.Range("H140").Value = 1421.7826
$ws.Range("I140").Value = 904.25
$ws.Range("J140").Value = 2604.7144
$ws.Range("K140").Value = 2712.75
$ws.Range("L140").Value = 7814.1432
$ws.Range("M140").Value = 2467.25
$ws.Range("N140").Value = -18174.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5517.391
$ws.Range("I80").Value = 5885
$ws.Range("K80").Value = 5885
$ws.Range("M80").Value = -4887
$ws.Range("H83").Value = 5517.391
$ws.Range("I83").Value = 5885
$ws.Range("K83").Value = 29425
$ws.Range("M83").Value = -24433

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1969
$ws.Range("I93").Value = 2017
$ws.Range("J93").Value = 1825
$ws.Range("K93").Value = 2017
$ws.Range("L93").Value = 1825
$ws.Range("M93").Value = -769
$ws.Range("N93").Value = -4321

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 23949.75
$ws.Range("J64").Value = 23949.75
$ws.Range("L64").Value = 23949.75
$ws.Range("N64").Value = -24445.75
$ws.Range("H67").Value = 23949.75
$ws.Range("J67").Value = 23949.75
$ws.Range("L67").Value = 23949.75
$ws.Range("N67").Value = -25665.75
$ws.Range("H132").Value = 28354.744
$ws.Range("I132").Value = 92864.37
$ws.Range("J132").Value = 3011.6785
$ws.Range("K132").Value = 278593.11
$ws.Range("L132").Value = 9035.0355
$ws.Range("M132").Value = -276063.11
$ws.Range("N132").Value = -14095.0355
